$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new Strategy text for the "Reverse a Linked List" row (row 8), in column D,
# matching the style (wrap text) already used in column D for other rows.
$ws.Range("D8").Value = "1. Use Iterative approach with 2 ref's prev and curr`n2. Use Recursive approach by reversing rest n-1 nodes and linking head node`n3. Use Recursive approach by reversing first n-1 nodes"
$ws.Range("D8").WrapText = $true

# Increase row height to fit the new multi-line content.
$ws.Rows.Item(8).RowHeight = 87

# Update the selection / view state to match the recorded edit location.
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("C8:C10").Select()
$excel.ActiveCell = $ws.Range("C10")
